# Auto-generated Excel COM-interop script
# Applies the "data up to 10th" update to fb-survey-communityState workbook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in previously-missing Puerto Rico (AR) values for rows 159-162 ---
$ws.Range("AR159").Value = 18.6363636
$ws.Range("AR160").Value = 11.832048
$ws.Range("AR161").Value = 11.3282489
$ws.Range("AR162").Value = 16.4

# --- Correct recomputed smoothed values in row 186 (W, AB, AX) ---
$ws.Range("W186").Value = 11.2812968
$ws.Range("AB186").Value = 25.4832589
$ws.Range("AX186").Value = 26.9105303

# --- New daily survey data rows: 04 Aug 2020 (187) through 10 Aug 2020 (193) ---

$row187 = New-Object 'object[,]' 1,56
$row187[0,0] = 20.3557312
$row187[0,1] = 35.7011163
$row187[0,2] = 28.3207753
$row187[0,3] = $null
$row187[0,4] = 26.6510193
$row187[0,5] = 20.2793995
$row187[0,6] = 18.5995839
$row187[0,7] = 9.389219499999999
$row187[0,8] = 13.4101942
$row187[0,9] = 17.1052632
$row187[0,10] = 27.2699425
$row187[0,11] = 32.5519121
$row187[0,12] = $null
$row187[0,13] = 14.2563198
$row187[0,14] = 28.5352863
$row187[0,15] = 31.8284205
$row187[0,16] = 19.4136951
$row187[0,17] = 24.5046699
$row187[0,18] = 26.2906227
$row187[0,19] = 25.4281233
$row187[0,20] = 34.8426969
$row187[0,21] = 11.2971668
$row187[0,22] = 15.5127651
$row187[0,23] = 9.913322600000001
$row187[0,24] = 16.6008978
$row187[0,25] = 20.7643292
$row187[0,26] = 25.3432738
$row187[0,27] = $null
$row187[0,28] = 37.2834541
$row187[0,29] = 24.1856456
$row187[0,30] = 20.7421065
$row187[0,31] = 25.1169563
$row187[0,32] = 25.3340167
$row187[0,33] = 10.1552043
$row187[0,34] = 11.8648329
$row187[0,35] = 18.7568802
$row187[0,36] = 23.8476356
$row187[0,37] = 11.6539355
$row187[0,38] = 21.7925562
$row187[0,39] = 29.5966904
$row187[0,40] = 14.9991643
$row187[0,41] = 14.366885
$row187[0,42] = $null
$row187[0,43] = 14.3772381
$row187[0,44] = 29.5169825
$row187[0,45] = 23.7142107
$row187[0,46] = 29.946227
$row187[0,47] = 30.8981091
$row187[0,48] = 27.1207642
$row187[0,49] = 18.6137295
$row187[0,50] = $null
$row187[0,51] = 8.638551100000001
$row187[0,52] = 17.0405936
$row187[0,53] = 21.2255787
$row187[0,54] = 21.5841355
$row187[0,55] = 21.8377106
$ws.Range("B187:BE187").Value = $row187

$row188 = New-Object 'object[,]' 1,56
$row188[0,0] = 20.9615385
$row188[0,1] = 35.9064994
$row188[0,2] = 28.7458818
$row188[0,3] = $null
$row188[0,4] = 26.4814908
$row188[0,5] = 20.0646674
$row188[0,6] = 18.5354695
$row188[0,7] = 9.342936099999999
$row188[0,8] = 13.3605601
$row188[0,9] = 17.1173763
$row188[0,10] = 27.0046123
$row188[0,11] = 32.1605344
$row188[0,12] = $null
$row188[0,13] = 14.4533486
$row188[0,14] = 28.214731
$row188[0,15] = 31.8266377
$row188[0,16] = 19.3933465
$row188[0,17] = 24.6730327
$row188[0,18] = 27.2088415
$row188[0,19] = 25.5517724
$row188[0,20] = 34.2205393
$row188[0,21] = 11.1666199
$row188[0,22] = 15.7333125
$row188[0,23] = 9.347982699999999
$row188[0,24] = 16.5564355
$row188[0,25] = 21.165553
$row188[0,26] = 25.7079996
$row188[0,27] = $null
$row188[0,28] = 36.8244932
$row188[0,29] = 25.5614527
$row188[0,30] = 20.8144657
$row188[0,31] = 25.6123352
$row188[0,32] = 26.5705539
$row188[0,33] = 9.7436864
$row188[0,34] = 11.8298294
$row188[0,35] = 18.3476935
$row188[0,36] = 24.543586
$row188[0,37] = 11.6203848
$row188[0,38] = 21.709507
$row188[0,39] = 29.8252549
$row188[0,40] = 14.919132
$row188[0,41] = 14.8106893
$row188[0,42] = $null
$row188[0,43] = 13.3874395
$row188[0,44] = 30.0487507
$row188[0,45] = 24.0628165
$row188[0,46] = 29.8433123
$row188[0,47] = 30.2114157
$row188[0,48] = 26.5512724
$row188[0,49] = 18.5047112
$row188[0,50] = $null
$row188[0,51] = 9.195845
$row188[0,52] = 17.3285811
$row188[0,53] = 20.7821335
$row188[0,54] = 22.6614474
$row188[0,55] = 20.5592813
$ws.Range("B188:BE188").Value = $row188

$ws.Range("A189").Value = "06 08 2020"
$row189 = New-Object 'object[,]' 1,56
$row189[0,0] = 20.2323718
$row189[0,1] = 34.8864049
$row189[0,2] = 29.1666673
$row189[0,3] = $null
$row189[0,4] = 26.2362149
$row189[0,5] = 19.9550069
$row189[0,6] = 17.9468251
$row189[0,7] = 9.917934300000001
$row189[0,8] = 13.7152778
$row189[0,9] = 15.9334461
$row189[0,10] = 26.8284194
$row189[0,11] = 32.4997761
$row189[0,12] = $null
$row189[0,13] = 14.8011364
$row189[0,14] = 28.3188429
$row189[0,15] = 30.9130974
$row189[0,16] = 19.6469514
$row189[0,17] = 25.1457936
$row189[0,18] = 26.4153
$row189[0,19] = 26.1302564
$row189[0,20] = 34.206333
$row189[0,21] = 10.7842124
$row189[0,22] = 15.4791993
$row189[0,23] = 10.4313871
$row189[0,24] = 16.7483337
$row189[0,25] = 20.5368202
$row189[0,26] = 25.6943563
$row189[0,27] = $null
$row189[0,28] = 36.1808014
$row189[0,29] = 24.4569456
$row189[0,30] = 20.6513023
$row189[0,31] = 26.9272471
$row189[0,32] = 27.4644938
$row189[0,33] = 9.996009600000001
$row189[0,34] = 11.9295588
$row189[0,35] = 19.1843197
$row189[0,36] = 24.7130041
$row189[0,37] = 11.492425
$row189[0,38] = 21.7834778
$row189[0,39] = 29.8271828
$row189[0,40] = 14.7566226
$row189[0,41] = 14.8538277
$row189[0,42] = $null
$row189[0,43] = 13.456802
$row189[0,44] = 29.9102766
$row189[0,45] = 25.3737107
$row189[0,46] = 29.8950326
$row189[0,47] = 30.0509898
$row189[0,48] = 26.3290567
$row189[0,49] = 18.5002105
$row189[0,50] = $null
$row189[0,51] = 8.263513
$row189[0,52] = 17.340269
$row189[0,53] = 20.7120269
$row189[0,54] = 22.2560424
$row189[0,55] = 21.9352936
$ws.Range("B189:BE189").Value = $row189

$ws.Range("A190").Value = "07 08 2020"
$row190 = New-Object 'object[,]' 1,56
$row190[0,0] = 20.081663
$row190[0,1] = 35.0305403
$row190[0,2] = 29.002053
$row190[0,3] = $null
$row190[0,4] = 26.0933777
$row190[0,5] = 19.7633532
$row190[0,6] = 17.5124234
$row190[0,7] = 9.914285700000001
$row190[0,8] = 13.6206897
$row190[0,9] = 15.4996826
$row190[0,10] = 26.4308234
$row190[0,11] = 32.1094825
$row190[0,12] = $null
$row190[0,13] = 14.8708081
$row190[0,14] = 27.7252229
$row190[0,15] = 30.4400575
$row190[0,16] = 19.7355959
$row190[0,17] = 24.5584956
$row190[0,18] = 25.7054229
$row190[0,19] = 26.1140484
$row190[0,20] = 34.0206967
$row190[0,21] = 10.8152434
$row190[0,22] = 15.6529988
$row190[0,23] = 10.4923184
$row190[0,24] = 16.7649258
$row190[0,25] = 20.92824
$row190[0,26] = 26.5155362
$row190[0,27] = $null
$row190[0,28] = 36.2971219
$row190[0,29] = 25.7843051
$row190[0,30] = 20.3181243
$row190[0,31] = 25.3499071
$row190[0,32] = 26.7977587
$row190[0,33] = 10.2930267
$row190[0,34] = 11.8354834
$row190[0,35] = 19.3668077
$row190[0,36] = 24.7010005
$row190[0,37] = 11.724194
$row190[0,38] = 21.7565851
$row190[0,39] = 30.2633993
$row190[0,40] = 14.822062
$row190[0,41] = 14.7155146
$row190[0,42] = $null
$row190[0,43] = 13.1936266
$row190[0,44] = 29.5462766
$row190[0,45] = 24.3514678
$row190[0,46] = 29.9090546
$row190[0,47] = 29.7899525
$row190[0,48] = 25.9627901
$row190[0,49] = 18.5381088
$row190[0,50] = $null
$row190[0,51] = 9.8120388
$row190[0,52] = 17.3004219
$row190[0,53] = 20.9679286
$row190[0,54] = 22.1188847
$row190[0,55] = 20.801157
$ws.Range("B190:BE190").Value = $row190

$ws.Range("A191").Value = "08 08 2020"
$row191 = New-Object 'object[,]' 1,56
$row191[0,0] = 20.0242718
$row191[0,1] = 34.815552
$row191[0,2] = 29.0144283
$row191[0,3] = $null
$row191[0,4] = 25.5482216
$row191[0,5] = 19.6736523
$row191[0,6] = 17.2241917
$row191[0,7] = 10.0829795
$row191[0,8] = 13.2738095
$row191[0,9] = 15.7110327
$row191[0,10] = 25.9399211
$row191[0,11] = 31.7708947
$row191[0,12] = $null
$row191[0,13] = 15.323026
$row191[0,14] = 27.4219389
$row191[0,15] = 30.993375
$row191[0,16] = 19.8870186
$row191[0,17] = 24.5077037
$row191[0,18] = 25.5518819
$row191[0,19] = 26.4600273
$row191[0,20] = 33.9481835
$row191[0,21] = 10.7680299
$row191[0,22] = 15.0109604
$row191[0,23] = 10.7091433
$row191[0,24] = 16.7964911
$row191[0,25] = 20.9797716
$row191[0,26] = 25.5447154
$row191[0,27] = $null
$row191[0,28] = 36.8428023
$row191[0,29] = 24.735524
$row191[0,30] = 20.8978711
$row191[0,31] = 26.3059082
$row191[0,32] = 27.5909316
$row191[0,33] = 10.5971897
$row191[0,34] = 11.8298617
$row191[0,35] = 18.8525915
$row191[0,36] = 24.9207334
$row191[0,37] = 11.7340101
$row191[0,38] = 21.7213358
$row191[0,39] = 29.3569052
$row191[0,40] = 14.8383017
$row191[0,41] = 14.8056069
$row191[0,42] = $null
$row191[0,43] = 13.6640075
$row191[0,44] = 29.4157197
$row191[0,45] = 23.3646901
$row191[0,46] = 29.0117651
$row191[0,47] = 29.5754224
$row191[0,48] = 26.2180531
$row191[0,49] = 18.634555
$row191[0,50] = $null
$row191[0,51] = 9.983718
$row191[0,52] = 17.2370206
$row191[0,53] = 20.9881611
$row191[0,54] = 20.9211204
$row191[0,55] = 20.353184
$ws.Range("B191:BE191").Value = $row191

$ws.Range("A192").Value = "09 08 2020"
$row192 = New-Object 'object[,]' 1,56
$row192[0,0] = 19.7665663
$row192[0,1] = 33.9954347
$row192[0,2] = 28.9760418
$row192[0,3] = $null
$row192[0,4] = 25.4486435
$row192[0,5] = 19.4227517
$row192[0,6] = 17.1530023
$row192[0,7] = 10.0289017
$row192[0,8] = 14.0022676
$row192[0,9] = 15.0575375
$row192[0,10] = 25.9505333
$row192[0,11] = 31.8989296
$row192[0,12] = $null
$row192[0,13] = 15.2233115
$row192[0,14] = 27.4019876
$row192[0,15] = 30.6362489
$row192[0,16] = 20.0534445
$row192[0,17] = 24.4245511
$row192[0,18] = 26.1036784
$row192[0,19] = 26.5472952
$row192[0,20] = 33.6057561
$row192[0,21] = 10.6418086
$row192[0,22] = 14.7220758
$row192[0,23] = 10.5008945
$row192[0,24] = 16.6583151
$row192[0,25] = 20.4308228
$row192[0,26] = 25.8887085
$row192[0,27] = $null
$row192[0,28] = 35.7021528
$row192[0,29] = 24.2662218
$row192[0,30] = 20.6905413
$row192[0,31] = 25.2791554
$row192[0,32] = 27.5583378
$row192[0,33] = 9.877212
$row192[0,34] = 11.9388498
$row192[0,35] = 19.2530937
$row192[0,36] = 25.0687272
$row192[0,37] = 11.2951314
$row192[0,38] = 21.4788104
$row192[0,39] = 29.5767771
$row192[0,40] = 15.2125601
$row192[0,41] = 14.8439081
$row192[0,42] = $null
$row192[0,43] = 13.2873839
$row192[0,44] = 28.4475278
$row192[0,45] = 23.1972295
$row192[0,46] = 29.2849651
$row192[0,47] = 29.0155968
$row192[0,48] = 25.7695746
$row192[0,49] = 18.4505202
$row192[0,50] = $null
$row192[0,51] = 9.4638692
$row192[0,52] = 17.0900662
$row192[0,53] = 20.8283298
$row192[0,54] = 21.648955
$row192[0,55] = 20.0499912
$ws.Range("B192:BE192").Value = $row192

# --- Final row: date label only (10 Aug 2020), no data yet ---
$ws.Range("A193").Value = "10 08 2020"
